# Commit: "added K_CSM and K_PSM results printing"
#
# - Remove the embedded feasibility-tensor picture from Sheet1 (it's being
#   superseded by printed K_CSM / K_PSM results, so the plot image + its
#   drawing anchor go away).
# - Flip the Flag column (D) from 0 -> 1 for the rows that are now
#   considered feasible once K_CSM/K_PSM are accounted for.
# - Leave the active selection on V5 (where the cursor was when the file
#   was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the embedded picture / any other shapes anchored on the sheet.
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Rows whose D (Flag) value flips from 0 to 1.
$rows = @(
    73,
    86, 87, 88, 89,
    91, 92, 93,
    95, 96, 97,
    99, 100, 101,
    104, 105,
    108, 109,
    112, 113, 114, 115, 116, 117, 118, 119, 120, 121,
    122, 123, 124, 125, 126, 127, 128, 129,
    131, 132, 133,
    135, 136, 137,
    139, 140, 141, 142, 143, 144, 145, 146, 147, 148, 149, 150,
    151, 152, 153, 154, 155, 156, 157, 158, 159, 160, 161, 162,
    163, 164, 165,
    167, 168, 169
)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Restore the saved selection/active cell.
$ws.Range("V5").Select()
